$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new column before column B. This shifts the existing
# description column (old B) to column C, carrying its width/bestFit
# formatting along with it, and leaves column A free for new data.
$ws.Range("B:B").Insert()

# Row 1: existing first entry, now tagged with a completion date and
# "Done" status.
$ws.Range("A1").Value = 45922
$ws.Range("A1").NumberFormat = "mm-dd-yy"
$ws.Range("A1").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("B1").Value = "Done"
$ws.Range("B1").Interior.Color = $ws.Range("C1").Interior.Color

# Row 2: new entry - "Change Requested Start Date..." marked done, with note.
$ws.Range("A2").Value = 45922
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("B2").Value = "Done"
$ws.Range("B2").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("C2").Value = "Change Requested Start Date to Include Today and up to 29 Days"
$ws.Range("C2").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("D2").Value = "per JoAnna and John keep it at 7 days"

# Row 3: new entry - "Add Enrollment to Transaction" marked done.
$ws.Range("A3").Value = 45922
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("B3").Value = "Done"
$ws.Range("B3").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("C3").Value = "Add Enrollment to Transaction"
$ws.Range("C3").Interior.Color = $ws.Range("C1").Interior.Color

# Row 8: new entry - "Monaco Remove Tennis Passport and Racquet Addons", with note.
$ws.Range("A8").Value = 45922
$ws.Range("A8").NumberFormat = "mm-dd-yy"
$ws.Range("A8").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("B8").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("C8").Interior.Color = $ws.Range("C1").Interior.Color
$ws.Range("D8").Value = "removed from inventory per John"

# Size the new date column to fit its contents, like the description
# column next to it.
$ws.Columns.Item(1).AutoFit()

$ws.Range("D11").Select()
